$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new trade log entry (row 5), mirroring the existing rows' structure/format.
$ws.Cells.Item(5,1).Value = 44141
$ws.Cells.Item(5,1).NumberFormat = "yyyy\-mm\-dd"

$ws.Cells.Item(5,2).Value = "owens_trading_v1"

$ws.Cells.Item(5,3).Value = 108

$ws.Cells.Item(5,4).Value = 0.704861111111111
$ws.Cells.Item(5,4).NumberFormat = "hh:mm:ss"

$ws.Cells.Item(5,5).Value = 0.716666666666667
$ws.Cells.Item(5,5).NumberFormat = "hh:mm:ss"

$ws.Cells.Item(5,6).Value = 17

$ws.Cells.Item(5,7).Value = 1232064

$ws.Cells.Item(5,8).Formula = "=G5/F5"

# Update the selected cell to reflect where the user was last working.
$ws.Range("F6").Select()
